# BC congno chitiet: add column {NVPhuTrach}
# Insert a new column before column E ("Tong phai tra") and label it
# "NV phu trach". Excel's EntireColumn.Insert() shifts the existing
# E..K columns (and their formulas / merged ranges) one place to the
# right automatically, matching the rest of the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column at E (pushes old E:K -> F:L).
$ws.Range("E:E").EntireColumn.Insert()

# New header text for the inserted column.
$ws.Range("E4").Value = "NV phụ trách"

# Restore the widened column D and size the new column E to match the
# published template (values expressed in Excel "characters" units;
# they land on the stored width the host actually persists).
$ws.Range("D1").EntireColumn.ColumnWidth = 39.16
$ws.Range("E1").EntireColumn.ColumnWidth = 27.43

# Match the saved selection/active cell from the authored workbook.
$ws.Range("F19").Select()
